# Revert "Remove aviation exemption from LCFS"
# Restores the full source citation on the About sheet, re-adds the
# "Based on the California LCFS, we choose to exempt aircraft." note,
# and exempts aircraft (sets BVTStL values to 0) on the data sheet.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("About")
$ws2 = $wb.Worksheets.Item("BVTStL")

# --- About sheet -----------------------------------------------------
# Make room for the expanded source citation: push the "Notes" block
# (currently rows 5-9) down to rows 9-13 by inserting 4 rows above it.
$ws1.Rows("5:8").Insert()

# Expand the "Source:" entry from the bare "None" placeholder into a
# full citation spread across B3:B7.
$ws1.Range("B3").Value = "California Air Resources Board"
$ws1.Range("B4").Value = 2015
$ws1.Range("B4").HorizontalAlignment = -4131
$ws1.Range("B5").Value = "Low Carbon Fuel Standard: Final Regulation Order"
$ws1.Range("B6").Value = "https://www.arb.ca.gov/regact/2015/lcfs2015/lcfsfinalregorder.pdf"
$ws1.Range("B7").Value = "Page 15"

# Re-add the explanatory note about exempting aircraft (row 14 stays
# blank, mirroring the original workbook's layout).
$ws1.Range("A15").Value = "Based on the California LCFS, we choose to exempt aircraft."

# --- BVTStL sheet ------------------------------------------------------
# Aircraft (row 4) are exempt from the LCFS: both passenger and freight
# flags go back to 0.
$ws2.Range("B4").Value = 0
$ws2.Range("C4").Value = 0
